{"js": "// 1) \"Occurrences / Objects Contexts.\" -> \"Occurrences / Objects / Contexts / SPOs.\"\nconst oldSnippet1 = \"Occurrences / Objects Contexts\";\nconst newSnippet1 = \"Occurrences / Objects / Contexts / SPOs\";\n\nconst results1 = context.document.body.search(oldSnippet1, { matchCase: true });\nresults1.load(\"text\");\nawait context.sync();\n\nif (results1.items.length > 0) {\n  results1.items[0].insertText(newSnippet1, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Append an extra sentence to the \"Role: occurrence / object in CSPO slots...\" paragraph.\nconst oldSnippet2 =\n  \"Role: occurrence / object in CSPO slots. Denotes resource types in positions in statements (i.e.: Kind in Relation).\";\nconst newSnippet2 =\n  \"Role: occurrence / object in CSPO slots. Denotes resource types in positions in statements (i.e.: Kind in Relation). Role CSPO is object / occurrence in statement occurrence position, Role type (i.e.: Kind, Relation) stated as Role instances in Meta Models with corresponding Kinds for its complimentary CSPO resources.\";\n\nconst results2 = context.document.body.search(oldSnippet2, { matchCase: true });\nresults2.load(\"text\");\nawait context.sync();\n\nif (results2.items.length > 0) {\n  results2.items[0].insertText(newSnippet2, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) \"Occurrences / Objects Contexts\" -> \"Occurrences / Objects / Contexts / SPOs\"\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Execute(\n  \"Occurrences / Objects Contexts\",\n  $false,\n  $false,\n  $false,\n  $false,\n  $false,\n  $true,\n  0,\n  $false,\n  \"Occurrences / Objects / Contexts / SPOs\",\n  2\n)\n\n# 2) Append a trailing sentence to the \"Role: occurrence / object in CSPO slots...\" paragraph.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Execute(\n  \"Role: occurrence / object in CSPO slots. Denotes resource types in positions in statements (i.e.: Kind in Relation).\",\n  $false,\n  $false,\n  $false,\n  $false,\n  $false,\n  $true,\n  0,\n  $false,\n  \"Role: occurrence / object in CSPO slots. Denotes resource types in positions in statements (i.e.: Kind in Relation). Role CSPO is object / occurrence in statement occurrence position, Role type (i.e.: Kind, Relation) stated as Role instances in Meta Models with corresponding Kinds for its complimentary CSPO resources.\",\n  2\n)\n"}
